# Order.xlsx update
# - Refresh order data rows 2-9 with new values
# - Add a new order row 10
# - Adjust sheet view (zoom) and selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "A123456B5678784(7)"
$ws.Cells.Item(2,2).Value = 222001
$ws.Cells.Item(2,3).Value = "1126 2, food2 1, food3 4, food4 2, food5 3, tools1 10, tools2 2, necessities1 3, item4 1, necessities2 4, item5 1"
$ws.Cells.Item(2,4).Value = "2022-02-10"
$ws.Cells.Item(2,5).Value = 987654
$ws.Cells.Item(2,6).Value = "Credit Card"
$ws.Cells.Item(2,7).Value = "In Transit"
$ws.Cells.Item(2,8).Value = "T"

# Row 3
$ws.Cells.Item(3,1).Value = "C235786A1204756(8)"
$ws.Cells.Item(3,2).Value = 222005
$ws.Cells.Item(3,3).Value = "tools3 1, necessities5 3, food2 2, food2 1, food3 5"
$ws.Cells.Item(3,4).Value = "2022-05-12"
$ws.Cells.Item(3,5).Value = 314778
$ws.Cells.Item(3,6).Value = "Mobile"
$ws.Cells.Item(3,7).Value = "In Transit"
$ws.Cells.Item(3,8).Value = "F"

# Row 4
$ws.Cells.Item(4,1).Value = "B250035C1578965(9)"
$ws.Cells.Item(4,2).Value = 222004
$ws.Cells.Item(4,3).Value = "necessities1 4, item1 10, necessities5 1"
$ws.Cells.Item(4,4).Value = "2021-10-07"
$ws.Cells.Item(4,5).Value = 412789
$ws.Cells.Item(4,6).Value = "FPS"
$ws.Cells.Item(4,7).Value = "Received"
$ws.Cells.Item(4,8).Value = "T"

# Row 5
$ws.Cells.Item(5,1).Value = "XR785214A2533687(5)"
$ws.Cells.Item(5,2).Value = 222006
$ws.Cells.Item(5,3).Value = "food1 2, food3 1, food3 7"
$ws.Cells.Item(5,4).Value = "2023-01-05"
$ws.Cells.Item(5,5).Value = 203578
$ws.Cells.Item(5,6).Value = "E Wallet"
$ws.Cells.Item(5,7).Value = "Received"
$ws.Cells.Item(5,8).Value = "T"

# Row 6
$ws.Cells.Item(6,1).Value = "K128746B2547896(1)"
$ws.Cells.Item(6,2).Value = 222303
$ws.Cells.Item(6,3).Value = "1133 5, 1139 1, food2 5"
$ws.Cells.Item(6,4).Value = "2020-04-17"
$ws.Cells.Item(6,5).Value = 155985
$ws.Cells.Item(6,6).Value = "Mobile"
$ws.Cells.Item(6,7).Value = "Received"
$ws.Cells.Item(6,8).Value = "F"

# Row 7
$ws.Cells.Item(7,1).Value = "TN325221B7895125(4)"
$ws.Cells.Item(7,2).Value = 222104
$ws.Cells.Item(7,3).Value = "food10 3"
$ws.Cells.Item(7,4).Value = "2023-02-15"
$ws.Cells.Item(7,5).Value = 984203
$ws.Cells.Item(7,6).Value = "FPS"
$ws.Cells.Item(7,7).Value = "Received"
$ws.Cells.Item(7,8).Value = "T"

# Row 8
$ws.Cells.Item(8,1).Value = "NN985234C5324800(3)"
$ws.Cells.Item(8,2).Value = 222103
$ws.Cells.Item(8,3).Value = "food11 28, necessities2 2, food2 1, food3 10, item4 5, item3 7, tool1 1"
$ws.Cells.Item(8,4).Value = "2021-05-08"
$ws.Cells.Item(8,5).Value = 741265
$ws.Cells.Item(8,6).Value = "E Wallet"
$ws.Cells.Item(8,7).Value = "In Transit"
$ws.Cells.Item(8,8).Value = "T"

# Row 9
$ws.Cells.Item(9,1).Value = "U129658A9357825(6)"
$ws.Cells.Item(9,2).Value = 222002
$ws.Cells.Item(9,3).Value = "item2 10, item3 1, food1 10, food2 5"
$ws.Cells.Item(9,4).Value = "2022-02-17"
$ws.Cells.Item(9,5).Value = 354154
$ws.Cells.Item(9,6).Value = "Mobile"
$ws.Cells.Item(9,7).Value = "In Transit"
$ws.Cells.Item(9,8).Value = "F"

# Row 10 (new)
$ws.Cells.Item(10,1).Value = "MN123254C8548952(1)"
$ws.Cells.Item(10,2).Value = 222001
$ws.Cells.Item(10,3).Value = "1123 1, 1131 2, tools5 5"
$ws.Cells.Item(10,4).Value = "2023-03-03"
$ws.Cells.Item(10,5).Value = 298732
$ws.Cells.Item(10,6).Value = "Mobile"
$ws.Cells.Item(10,7).Value = "Received"
$ws.Cells.Item(10,8).Value = "T"

# View adjustments: zoom to 94% and move selection back to A1
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 94
